$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B ("ASIN") to host "Week_Start_Date".
$ws1.Columns.Item(2).Insert()

# New header
$ws1.Range("B1").Value = "Week_Start_Date"

# Week labels without the leading zero, e.g. "W01" -> "W1"
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")

# Week start dates (kept as plain text, not Excel dates)
$weekStartDates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)

# Updated "MyForecast" values (now in column D after the column insert)
$myForecast = @(687,682,642,620,632,649,653,658,643,638,636,636,639,645,653,666)

for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $weeks[$i]
    # Leading apostrophe forces text storage instead of an Excel date serial.
    $ws1.Cells.Item($r, 2).Value = "'" + $weekStartDates[$i]
    $ws1.Cells.Item($r, 4).Value = $myForecast[$i]
    $ws1.Cells.Item($r, 10).Value = $false
}

# ---------------------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value  = "'10380"
$ws2.Range("B10").Value = "'5223"
$ws2.Range("B11").Value = "'2632"
$ws2.Range("B12").Value = "'687"
$ws2.Range("B14").Value = "'620"
